$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the cryptos list update: price/volume refresh and a new "LEO" row
# inserted at row 28 (pushing subsequent rows down and dropping the last
# row, SuiNetwork, off the bottom of the A1:E51 range).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.886.13"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.502.05"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.44"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.83"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.498.78"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.56"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.98%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.092.68"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.906.84"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.491.98"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.49"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.38"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.77%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.37"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.09"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.641.78"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.18%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.89%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.69"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.99%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.96"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.60"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.50"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.64"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.171"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.60"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.14"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.85"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.497.03"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.98"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.33"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.67%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "176.76"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0898"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.41"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.895"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "30.37"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.72%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.67"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.28"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.53"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.82%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.61"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.08%  "
